$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 ("Marking") updates
$ws.Range("B11").Value = 9
$ws.Range("C11").Value = 2

# Row 12 ("Total") updates
$ws.Range("B12").Value = 234
$ws.Range("C12").Value = -4

# E12 holds a text fraction-like string, force text so Excel does not
# reinterpret "230/252" as a date value.
$ws.Range("E12").Value = "230/252"
